$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2027.678
$ws.Range("J17").Value = 2045.3966
$ws.Range("L17").Value = 6136.1898
$ws.Range("N17").Value = -6472.1898

$ws.Range("H70").Value = 1515.5333
$ws.Range("I70").Value = 477.5
$ws.Range("J70").Value = 1675.2307
$ws.Range("K70").Value = 1432.5
$ws.Range("L70").Value = 5025.6921
$ws.Range("M70").Value = -1162.5
$ws.Range("N70").Value = -5565.6921

$ws.Range("H73").Value = 1515.5333
$ws.Range("I73").Value = 477.5
$ws.Range("J73").Value = 1675.2307
$ws.Range("K73").Value = 1432.5
$ws.Range("L73").Value = 5025.6921
$ws.Range("M73").Value = -496.5
$ws.Range("N73").Value = -6897.6921

$ws.Range("H74").Value = 4704.6665
$ws.Range("I74").Value = 4338
$ws.Range("K74").Value = 4338
$ws.Range("M74").Value = -3402

$ws.Range("H76").Value = 3566.6667
$ws.Range("I76").Value = 3342.8572
$ws.Range("K76").Value = 3342.8572
$ws.Range("M76").Value = -3027.8572

$ws.Range("H77").Value = 4704.6665
$ws.Range("I77").Value = 4338
$ws.Range("K77").Value = 21690
$ws.Range("M77").Value = -17010

$ws.Range("H79").Value = 3566.6667
$ws.Range("I79").Value = 3342.8572
$ws.Range("K79").Value = 3342.8572
$ws.Range("M79").Value = -2250.8572

$ws.Range("H98").Value = 1167.9166
$ws.Range("I98").Value = 1167.9166
$ws.Range("K98").Value = 1167.9166
$ws.Range("M98").Value = 330.0834

$ws.Range("H122").Value = 1167.9166
$ws.Range("I122").Value = 1167.9166
$ws.Range("K122").Value = 3503.7498
$ws.Range("M122").Value = -1053.7498

$ws.Range("H123").Value = 28780
$ws.Range("J123").Value = 28780
$ws.Range("L123").Value = 28780
$ws.Range("N123").Value = -38580

$ws.Range("H130").Value = 48545.59
$ws.Range("J130").Value = 48545.59
$ws.Range("L130").Value = 48545.59
$ws.Range("N130").Value = -58585.59

$ws.Range("H138").Value = 3827.5588
$ws.Range("I138").Value = 1802.4584
$ws.Range("J138").Value = 8687.799999999999
$ws.Range("K138").Value = 5407.3752
$ws.Range("L138").Value = 26063.4
$ws.Range("M138").Value = -267.3752000000004
$ws.Range("N138").Value = -36343.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 38118.332
$ws.Range("J24").Value = 38118.332
$ws.Range("L24").Value = 38118.332
$ws.Range("N24").Value = -38866.332

$ws.Range("H61").Value = 7332.5713
$ws.Range("J61").Value = 19292.1
$ws.Range("L61").Value = 19292.1
$ws.Range("N61").Value = -19716.1

$ws.Range("H100").Value = 38118.332
$ws.Range("J100").Value = 38118.332
$ws.Range("L100").Value = 38118.332
$ws.Range("N100").Value = -40282.332

$ws.Range("H136").Value = 7332.5713
$ws.Range("J136").Value = 19292.1
$ws.Range("L136").Value = 57876.3
$ws.Range("N136").Value = -62976.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1119142.2
$ws.Range("I105").Value = 1490655.2
$ws.Range("J105").Value = 4603.143
$ws.Range("K105").Value = 1490655.2
$ws.Range("L105").Value = 4603.143
$ws.Range("M105").Value = -1488908.2
$ws.Range("N105").Value = -8097.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1213743.2
$ws.Range("I58").Value = 1624451
$ws.Range("J58").Value = 3236.5264
$ws.Range("K58").Value = 1624451
$ws.Range("L58").Value = 3236.5264
$ws.Range("M58").Value = -1624248
$ws.Range("N58").Value = -3642.5264

$ws.Range("H132").Value = 3117.4385
$ws.Range("I132").Value = 3264.9
$ws.Range("J132").Value = 2796.8696
$ws.Range("K132").Value = 9794.700000000001
$ws.Range("L132").Value = 8390.6088
$ws.Range("M132").Value = -7264.700000000001
$ws.Range("N132").Value = -13450.6088

$ws.Range("H134").Value = 2373.9688
$ws.Range("I134").Value = 1364.475
$ws.Range("J134").Value = 4056.4583
$ws.Range("K134").Value = 4093.425
$ws.Range("L134").Value = 12169.3749
$ws.Range("M134").Value = -1558.425
$ws.Range("N134").Value = -17239.3749

$ws.Range("H136").Value = 1213743.2
$ws.Range("I136").Value = 1624451
$ws.Range("J136").Value = 3236.5264
$ws.Range("K136").Value = 4873353
$ws.Range("L136").Value = 9709.5792
$ws.Range("M136").Value = -4870803
$ws.Range("N136").Value = -14809.5792

$ws.Range("H141").Value = 24621.238
$ws.Range("I141").Value = 20263.111
$ws.Range("J141").Value = 27889.834
$ws.Range("K141").Value = 20263.111
$ws.Range("L141").Value = 27889.834
$ws.Range("M141").Value = -15083.111
$ws.Range("N141").Value = -38249.834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 411.17648
$ws.Range("I18").Value = 343.125
$ws.Range("K18").Value = 1029.375
$ws.Range("M18").Value = -860.375

$ws.Range("H102").Value = 6097.5454
$ws.Range("I102").Value = 5026
$ws.Range("J102").Value = 6499.375
$ws.Range("K102").Value = 15078
$ws.Range("L102").Value = 19498.125
$ws.Range("M102").Value = -12644
$ws.Range("N102").Value = -24366.125

$ws.Range("H107").Value = 1308.24
$ws.Range("I107").Value = 332.54544
$ws.Range("J107").Value = 2074.8572
$ws.Range("K107").Value = 997.63632
$ws.Range("L107").Value = 6224.571599999999
$ws.Range("M107").Value = 922.36368
$ws.Range("N107").Value = -10064.5716

$ws.Range("H122").Value = 697.875
$ws.Range("I122").Value = 441.07693
$ws.Range("K122").Value = 3969.69237
$ws.Range("M122").Value = -1519.69237

$ws.Range("H131").Value = 1678.4482
$ws.Range("J131").Value = 1498.3334
$ws.Range("L131").Value = 4495.0002
$ws.Range("N131").Value = -14575.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5927.44
$ws.Range("I70").Value = 5514.88
$ws.Range("J70").Value = 6340
$ws.Range("K70").Value = 5514.88
$ws.Range("L70").Value = 6340
$ws.Range("M70").Value = -5244.88
$ws.Range("N70").Value = -6880

$ws.Range("H73").Value = 5927.44
$ws.Range("I73").Value = 5514.88
$ws.Range("J73").Value = 6340
$ws.Range("K73").Value = 5514.88
$ws.Range("L73").Value = 6340
$ws.Range("M73").Value = -4578.88
$ws.Range("N73").Value = -8212

$ws.Range("H80").Value = 3875
$ws.Range("J80").Value = 3928.5715
$ws.Range("L80").Value = 3928.5715
$ws.Range("N80").Value = -5924.5715

$ws.Range("H83").Value = 3875
$ws.Range("J83").Value = 3928.5715
$ws.Range("L83").Value = 19642.8575
$ws.Range("N83").Value = -29626.8575

$ws.Range("H132").Value = 3643.7754
$ws.Range("I132").Value = 1670.5349
$ws.Range("J132").Value = 17785.334
$ws.Range("K132").Value = 5011.6047
$ws.Range("L132").Value = 53356.00199999999
$ws.Range("M132").Value = -2481.6047
$ws.Range("N132").Value = -58416.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6578.2188
$ws.Range("I7").Value = 4968.3687
$ws.Range("K7").Value = 4968.3687
$ws.Range("M7").Value = -4856.3687

$ws.Range("H16").Value = 1473
$ws.Range("I16").Value = 982.2143
$ws.Range("J16").Value = 2454.5715
$ws.Range("K16").Value = 982.2143
$ws.Range("L16").Value = 2454.5715
$ws.Range("M16").Value = -812.2143
$ws.Range("N16").Value = -2794.5715

$ws.Range("H32").Value = 10400
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 19800
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 19800
$ws.Range("M32").Value = -683
$ws.Range("N32").Value = -20434

$ws.Range("H126").Value = 6578.2188
$ws.Range("I126").Value = 4968.3687
$ws.Range("K126").Value = 14905.1061
$ws.Range("M126").Value = -12435.1061

$ws.Range("H132").Value = 4035.311
$ws.Range("I132").Value = 3917.1707
$ws.Range("J132").Value = 5246.25
$ws.Range("K132").Value = 11751.5121
$ws.Range("L132").Value = 15738.75
$ws.Range("M132").Value = -9221.5121
$ws.Range("N132").Value = -20798.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 20003894
$ws.Range("J14").Value = 4822
$ws.Range("L14").Value = 4822
$ws.Range("N14").Value = -5158

$ws.Range("H64").Value = 40114
$ws.Range("J64").Value = 40114
$ws.Range("L64").Value = 40114
$ws.Range("N64").Value = -40610

$ws.Range("H67").Value = 40114
$ws.Range("J67").Value = 40114
$ws.Range("L67").Value = 40114
$ws.Range("N67").Value = -41830

$ws.Range("H92").Value = 34966.668
$ws.Range("J92").Value = 34450
$ws.Range("L92").Value = 34450
$ws.Range("N92").Value = -39442

$ws.Range("H132").Value = 1674.7106
$ws.Range("I132").Value = 769.64
$ws.Range("J132").Value = 3415.2307
$ws.Range("K132").Value = 2308.92
$ws.Range("L132").Value = 10245.6921
$ws.Range("M132").Value = 221.0799999999999
$ws.Range("N132").Value = -15305.6921
